# Refresh the crypto price/volume table (GitHub Actions-style data pull).
# Updates the "Price" (D) and "Volume(1h)" (E) columns for most rows, and
# swaps the Quant/Maker rows (46-47) including their Coin name and Link.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Forces the value to be stored as text (matching the original inlineStr
    # cell type) even for strings that look numeric (e.g. "1.000", "0.9990"),
    # then resets the cell style so no stray NumberFormat sticks around.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '29.777.67'
$ws.Range("E2").Value = '  -0.52%  '
Set-TextValue "D3" '1.890.82'
$ws.Range("E3").Value = '  +0.12%  '
Set-TextValue "D4" '0.9990'
$ws.Range("E4").Value = '  -0.27%  '
Set-TextValue "D5" '0.7856'
$ws.Range("E5").Value = '  -4.97%  '
Set-TextValue "D6" '242.44'
$ws.Range("E6").Value = '  +0.62%  '
$ws.Range("E7").Value = '  -0.22%  '
Set-TextValue "D8" '0.3169'
$ws.Range("E8").Value = '  -1.25%  '
Set-TextValue "D9" '25.36'
$ws.Range("E9").Value = '  -4.31%  '
Set-TextValue "D10" '0.07024'
$ws.Range("E10").Value = '  +0.40%  '
Set-TextValue "D11" '0.08045'
$ws.Range("E11").Value = '  +0.18%  '
Set-TextValue "D12" '0.7655'
Set-TextValue "D13" '1.881.77'
$ws.Range("E13").Value = '  -1.67%  '
Set-TextValue "D14" '5.279'
$ws.Range("E14").Value = '  +1.92%  '
Set-TextValue "D15" '92.06'
$ws.Range("E15").Value = '  +0.07%  '
Set-TextValue "D16" '29.779.66'
$ws.Range("E16").Value = '  -0.54%  '
Set-TextValue "D17" '13.84'
$ws.Range("E17").Value = '  -1.01%  '
Set-TextValue "D18" '5.906'
$ws.Range("E18").Value = '  +0.23%  '
Set-TextValue "D19" '243.20'
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("E20").Value = '  -0.46%  '
Set-TextValue "D21" '8.195'
$ws.Range("E21").Value = '  +18.91%  '
$ws.Range("E22").Value = '  -0.21%  '
Set-TextValue "D23" '2.144.90'
$ws.Range("E23").Value = '  -0.63%  '
Set-TextValue "D24" '1.000'
$ws.Range("E24").Value = '  -0.17%  '
Set-TextValue "D25" '0.1651'
$ws.Range("E25").Value = '  +5.04%  '
Set-TextValue "D26" '9.295'
$ws.Range("E26").Value = '  +1.48%  '
Set-TextValue "D27" '165.71'
$ws.Range("E27").Value = '  -1.09%  '
Set-TextValue "D28" '18.67'
$ws.Range("E28").Value = '  -0.65%  '
$ws.Range("E29").Value = '  -1.55%  '
Set-TextValue "D30" '1.394'
$ws.Range("E30").Value = '  +1.33%  '
Set-TextValue "D31" '1.535'
$ws.Range("E31").Value = '  +1.16%  '
Set-TextValue "D32" '4.434'
$ws.Range("E32").Value = '  +4.73%  '
Set-TextValue "D33" '0.05611'
$ws.Range("E33").Value = '  -0.68%  '
Set-TextValue "D34" '4.034'
$ws.Range("E34").Value = '  -0.63%  '
$ws.Range("E35").Value = '  -0.61%  '
Set-TextValue "D36" '0.7361'
$ws.Range("E36").Value = '  +1.01%  '
$ws.Range("E37").Value = '  +0.22%  '
Set-TextValue "D38" '2.642'
$ws.Range("E38").Value = '  -2.92%  '
Set-TextValue "D39" '0.01906'
$ws.Range("E39").Value = '  +0.26%  '
Set-TextValue "D40" '2.769'
$ws.Range("E40").Value = '  +0.06%  '
Set-TextValue "D41" '0.4398'
$ws.Range("E41").Value = '  +0.20%  '
Set-TextValue "D42" '72.27'
$ws.Range("E42").Value = '  +0.92%  '
Set-TextValue "D43" '5.802'
$ws.Range("E43").Value = '  -2.10%  '
Set-TextValue "D44" '0.9992'
Set-TextValue "D45" '0.8381'
$ws.Range("E45").Value = '  -0.87%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D46" '1.022.45'
$ws.Range("E46").Value = '  +3.80%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D47" '102.18'
$ws.Range("E47").Value = '  +1.05%  '
$ws.Range("E48").Value = '  -1.15%  '
$ws.Range("E49").Value = '  +1.48%  '
Set-TextValue "D50" '7.407'
$ws.Range("E50").Value = '  -2.08%  '
Set-TextValue "D51" '2.037.39'
$ws.Range("E51").Value = '  -0.75%  '
